$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TESTS_WS")

# Row-level updates: column H = "Last Execution Started" (date serial),
# column I = "Execution time" (string like "4.12s").
$updates = @(
    @{ Row = 2; H = 42832.434130532405; I = "0.922s" },
    @{ Row = 3; H = 42832.43415622685; I = "10.005s" },
    @{ Row = 4; H = 42832.43428415509; I = "2.273s" },
    @{ Row = 8; H = 42832.43431512731; I = "12.539s" },
    @{ Row = 10; H = 42832.43446560185; I = "2.222s" },
    @{ Row = 12; H = 42832.43449332176; I = "6.618s" },
    @{ Row = 13; H = 42832.434574756946; I = "4.124s" },
    @{ Row = 14; H = 42832.43462458334; I = "2.176s" },
    @{ Row = 16; H = 42832.434659583334; I = "4.102s" },
    @{ Row = 17; H = 42832.43470921296; I = "4.12s" },
    @{ Row = 18; H = 42832.43475880787; I = "4.086s" },
    @{ Row = 21; H = 42832.43480920139; I = "4.097s" },
    @{ Row = 26; H = 42832.43486099537; I = "4.116s" },
    @{ Row = 28; H = 42832.434911412034; I = "4.106s" },
    @{ Row = 32; H = 42832.43496068287; I = "4.101s" },
    @{ Row = 33; H = 42832.43501011574; I = "4.098s" },
    @{ Row = 36; H = 42832.43505960648; I = "4.083s" },
    @{ Row = 43; H = 42832.435108819445; I = "4.107s" },
    @{ Row = 47; H = 42832.43515849537; I = "4.16s" },
    @{ Row = 48; H = 42832.43521607639; I = "4.211s" },
    @{ Row = 51; H = 42832.43526974537; I = "4.099s" },
    @{ Row = 54; H = 42832.43531903935; I = "4.077s" },
    @{ Row = 56; H = 42832.43536873843; I = "4.093s" },
    @{ Row = 58; H = 42832.43541854167; I = "4.105s" },
    @{ Row = 59; H = 42832.43546923611; I = "4.075s" },
    @{ Row = 61; H = 42832.43551818287; I = "2.23s" },
    @{ Row = 64; H = 42832.43554694444; I = "4.09s" },
    @{ Row = 66; H = 42832.43559803241; I = "4.151s" },
    @{ Row = 67; H = 42832.43564829861; I = "4.075s" },
    @{ Row = 68; H = 42832.435697106484; I = "4.06s" },
    @{ Row = 70; H = 42832.43574581019; I = "4.108s" },
    @{ Row = 72; H = 42832.435795590274; I = "4.113s" },
    @{ Row = 73; H = 42832.43584946759; I = "4.058s" },
    @{ Row = 74; H = 42832.43589824074; I = "4.064s" },
    @{ Row = 76; H = 42832.4359484838; I = "4.109s" },
    @{ Row = 77; H = 42832.43599893519; I = "4.119s" },
    @{ Row = 78; H = 42832.43604814815; I = "4.054s" },
    @{ Row = 79; H = 42832.43609695602; I = "4.103s" },
    @{ Row = 81; H = 42832.43614712963; I = "4.221s" },
    @{ Row = 82; H = 42832.43619803241; I = "4.092s" },
    @{ Row = 83; H = 42832.43624768518; I = "4.073s" },
    @{ Row = 85; H = 42832.43629765046; I = "4.082s" },
    @{ Row = 88; H = 42832.436346875; I = "4.094s" },
    @{ Row = 89; H = 42832.43639631944; I = "4.078s" },
    @{ Row = 91; H = 42832.43644554398; I = "31.849s" },
    @{ Row = 92; H = 42832.43681990741; I = "4.345s" },
    @{ Row = 93; H = 42832.43687329861; I = "0.034s" },
    @{ Row = 95; H = 42832.436874872685; I = "3.224s" },
    @{ Row = 96; H = 42832.43691486111; I = "4.265s" },
    @{ Row = 97; H = 42832.436968912036; I = "4.313s" },
    @{ Row = 100; H = 42832.43702385417; I = "0.159s" },
    @{ Row = 101; H = 42832.43702759259; I = "18.569s" },
    @{ Row = 104; H = 42832.43725040509; I = "22.337s" },
    @{ Row = 105; H = 42832.437516412036; I = "4.115s" },
    @{ Row = 107; H = 42832.437565810185; I = "2.138s" },
    @{ Row = 108; H = 42832.43759440972; I = "2.333s" },
    @{ Row = 110; H = 42832.43762435185; I = "7.145s" },
    @{ Row = 113; H = 42832.43771103009; I = "0.144s" },
    @{ Row = 118; H = 42832.437713796295; I = "4.149s" },
    @{ Row = 127; H = 42837.634543229164; I = "1.933s" }
)

foreach ($u in $updates) {
    $row = $u.Row
    $ws.Cells.Item($row, 8).Value = $u.H
    $ws.Cells.Item($row, 9).Value = $u.I
}

# Row 127 also got a different test suite / test case assigned.
$ws.Range("B127").Value = "EDELIVERY-2082 - Error when trying to download a message with an empty payload (disabled)"
$ws.Range("D127").Value = "Dom127-Submit Message-Basic Flow-Message with ID"
